# Update the "Date Placeholder" field text from 04-07-2023 to 11-07-2023
# across the slide master and every slide layout (mirrors what PowerPoint
# does when you edit the date shown via Insert > Header & Footer / by
# editing the Date placeholder directly on the master/layouts).

$p = $ppt.ActivePresentation

$oldDate = "04-07-2023"
$newDate = "11-07-2023"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master's own Date Placeholder shape.
Update-DatePlaceholder $master.Shapes

# Every slide layout owned by the master.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}
